$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two labels in column A
$ws.Range("A4").Value = "chosen_optimization_setting"
$ws.Range("A7").Value = "path_projects"

# Column B values (rows 2-9) are no longer populated; clear them
$ws.Range("B2:B9").ClearContents()

# Append three new label rows in column A
$ws.Range("A10").Value = "optimize_or_visualize"
$ws.Range("A11").Value = "chosen_visualization_setting"
$ws.Range("A12").Value = "optimization_or_visualization"

# Match the existing label formatting (border/bold/alignment) used by A1:A9
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the saved selection/active cell recorded in the workbook
$ws.Range("I14").Select() | Out-Null
